# Rename the three header/footer logo pictures:
#   - First-page footer's Pearson logo:  image2.png -> image1.png  (docPr id="3")
#   - Default footer's Pearson logo:     image2.png -> image1.png  (docPr id="2")
#   - First-page header's BTEC logo:     image1.jpg -> image2.jpg  (docPr id="1")
#
# wdHeaderFooterPrimary (default) = 1, wdHeaderFooterFirstPage = 2

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- First-page footer (contains "Authorised by: Head of BTEC Assessment") ---
$firstFooter = $sec.Footers.Item(2)
if ($firstFooter.Exists -and $firstFooter.Range.InlineShapes.Count -gt 0) {
    $firstFooter.Range.InlineShapes.Item(1).Name = "image1.png"
}

# --- Default footer (contains "Approved by: Delivery Manager") ---
$defaultFooter = $sec.Footers.Item(1)
if ($defaultFooter.Exists -and $defaultFooter.Range.InlineShapes.Count -gt 0) {
    $defaultFooter.Range.InlineShapes.Item(1).Name = "image1.png"
}

# --- First-page header (BTEC logo) ---
$firstHeader = $sec.Headers.Item(2)
if ($firstHeader.Exists -and $firstHeader.Range.InlineShapes.Count -gt 0) {
    $firstHeader.Range.InlineShapes.Item(1).Name = "image2.jpg"
}
